$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.3072063324736731
$ws.Range("C2").Value = 0.5773635797369776
$ws.Range("D2").Value = 0.6188915608581815
$ws.Range("E2").Value = 0.7866966129698166
$ws.Range("F2").Value = 0.733013573663008
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = 0.06855553159162397
$ws.Range("C3").Value = 0.8104126323077684
$ws.Range("D3").Value = 1.411859963597413
$ws.Range("E3").Value = 1.188217136552664
$ws.Range("F3").Value = 1.200974231718875
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = -0.2062133925061186
$ws.Range("C4").Value = 0.5585097485030701
$ws.Range("D4").Value = 0.5025698948149058
$ws.Range("E4").Value = 0.7089216422249399
$ws.Range("F4").Value = 0.6869075491602615
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = -0.02557072480903762
$ws.Range("C5").Value = 0.6736167808661312
$ws.Range("D5").Value = 0.868848526268419
$ws.Range("E5").Value = 0.9321204462237801
$ws.Range("F5").Value = 0.9439501535086928
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = -0.1229089188492291
$ws.Range("C6").Value = 0.4460013948249449
$ws.Range("D6").Value = 0.350767677791658
$ws.Range("E6").Value = 0.592256429084276
$ws.Range("F6").Value = 0.5871396736870387
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.05574729004729425
$ws.Range("C7").Value = 0.6872371239517284
$ws.Range("D7").Value = 0.9195149290446228
$ws.Range("E7").Value = 0.9589134106083942
$ws.Range("F7").Value = 0.9704962253316772
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = -0.1214114959447097
$ws.Range("C8").Value = 0.4048482891640412
$ws.Range("D8").Value = 0.2561152986148548
$ws.Range("E8").Value = 0.5060783522487944
$ws.Range("F8").Value = 0.4982679629511933
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.1522458912649086
$ws.Range("C9").Value = 0.6577382853197382
$ws.Range("D9").Value = 0.8416789472939559
$ws.Range("E9").Value = 0.9174306226053041
$ws.Range("F9").Value = 0.9179181168793588
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = -0.03691504020070672
$ws.Range("C10").Value = 0.3541807863119035
$ws.Range("D10").Value = 0.2034857836490605
$ws.Range("E10").Value = 0.4510939853833794
$ws.Range("F10").Value = 0.4563419822598951
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.2613592320042264
$ws.Range("C11").Value = 0.6789426103936371
$ws.Range("D11").Value = 1.020165729568819
$ws.Range("E11").Value = 1.010032538866357
$ws.Range("F11").Value = 0.9907586059223498
$ws.Range("G11").Value = 33
